# The experiment's stimulus images were switched from .jpg to .png files.
# Update every filename stored in the "Filename_Left" (D) / "Filename_Right" (E)
# columns of the orders sheet so they reference the new .png assets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("D", "E")) {
        $cell = $ws.Range("$col$r")
        $value = $cell.Value2
        if ($value -ne $null -and $value -is [string] -and $value.EndsWith(".jpg")) {
            $newValue = $value.Substring(0, $value.Length - 4) + ".png"
            $cell.Value = $newValue
        }
    }
}

# Reflect the saved selection state: the whole table is selected (anchored at A1)
# instead of the previous single-cell selection at K180.
$ws.Range("A1:F206").Select()
